$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "ECs"
$ws.Range("G2").Value = 5.441829000000001
$ws.Range("H2").Value = 16.325487
$ws.Range("I2").Value = 0.5729403216841985
$ws.Range("J2").Value = 0.5729403216841985
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.124002
$ws.Range("N2").Value = 0.372006
$ws.Range("O2").Value = 0.01161425268110074
$ws.Range("P2").Value = 0.01161425268110074
$ws.Range("Q2").Value = 0.6747976796580002
$ws.Range("R2").Value = 6.073179116922001
$ws.Range("S2").Value = 0.006654273667231425
$ws.Range("T2").Value = 0.006654273667231425

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("G3").Value = 5.441829000000001
$ws.Range("H3").Value = 16.325487
$ws.Range("I3").Value = 0.5729403216841985
$ws.Range("J3").Value = 0.5729403216841985
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.550379
$ws.Range("N3").Value = 16.651137
$ws.Range("O3").Value = 0.5198585843927942
$ws.Range("P3").Value = 0.5198585843927942
$ws.Range("Q3").Value = 30.20421340319101
$ws.Range("R3").Value = 271.837920628719
$ws.Range("S3").Value = 0.2978479445722996
$ws.Range("T3").Value = 0.2978479445722996

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 5.441829000000001
$ws.Range("H4").Value = 16.325487
$ws.Range("I4").Value = 0.5729403216841985
$ws.Range("J4").Value = 0.5729403216841985
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.002328333333334
$ws.Range("N4").Value = 15.006985
$ws.Range("O4").Value = 0.4685271629261051
$ws.Range("P4").Value = 0.4685271629261051
$ws.Range("Q4").Value = 27.22181539185501
$ws.Range("R4").Value = 244.996338526695
$ws.Range("S4").Value = 0.2684381034446676
$ws.Range("T4").Value = 0.2684381034446676

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "ECs"
$ws.Range("G5").Value = 2.579868
$ws.Range("H5").Value = 7.739604
$ws.Range("I5").Value = 0.2716201486343598
$ws.Range("J5").Value = 0.2716201486343598
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.124002
$ws.Range("N5").Value = 0.372006
$ws.Range("O5").Value = 0.01161425268110074
$ws.Range("P5").Value = 0.01161425268110074
$ws.Range("Q5").Value = 0.319908791736
$ws.Range("R5").Value = 2.879179125624
$ws.Range("S5").Value = 0.003154665039517596
$ws.Range("T5").Value = 0.003154665039517596

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("G6").Value = 2.579868
$ws.Range("H6").Value = 7.739604
$ws.Range("I6").Value = 0.2716201486343598
$ws.Range("J6").Value = 0.2716201486343598
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 5.550379
$ws.Range("N6").Value = 16.651137
$ws.Range("O6").Value = 0.5198585843927942
$ws.Range("P6").Value = 0.5198585843927942
$ws.Range("Q6").Value = 14.319245169972
$ws.Range("R6").Value = 128.873206529748
$ws.Range("S6").Value = 0.1412040659616186
$ws.Range("T6").Value = 0.1412040659616186

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("G7").Value = 2.579868
$ws.Range("H7").Value = 7.739604
$ws.Range("I7").Value = 0.2716201486343598
$ws.Range("J7").Value = 0.2716201486343598
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.002328333333334
$ws.Range("N7").Value = 15.006985
$ws.Range("O7").Value = 0.4685271629261051
$ws.Range("P7").Value = 0.4685271629261051
$ws.Range("Q7").Value = 12.90534679266
$ws.Range("R7").Value = 116.14812113394
$ws.Range("S7").Value = 0.1272614176332236
$ws.Range("T7").Value = 0.1272614176332236

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("D8").Value = "ECs"
$ws.Range("G8").Value = 1.476376
$ws.Range("H8").Value = 4.429128
$ws.Range("I8").Value = 0.1554395296814417
$ws.Range("J8").Value = 0.1554395296814417
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.124002
$ws.Range("N8").Value = 0.372006
$ws.Range("O8").Value = 0.01161425268110074
$ws.Range("P8").Value = 0.01161425268110074
$ws.Range("Q8").Value = 0.183073576752
$ws.Range("R8").Value = 1.647662190768
$ws.Range("S8").Value = 0.001805313974351723
$ws.Range("T8").Value = 0.001805313974351723

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("D9").Value = "FAPs"
$ws.Range("G9").Value = 1.476376
$ws.Range("H9").Value = 4.429128
$ws.Range("I9").Value = 0.1554395296814417
$ws.Range("J9").Value = 0.1554395296814417
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 5.550379
$ws.Range("N9").Value = 16.651137
$ws.Range("O9").Value = 0.5198585843927942
$ws.Range("P9").Value = 0.5198585843927942
$ws.Range("Q9").Value = 8.194446346504
$ws.Range("R9").Value = 73.75001711853599
$ws.Range("S9").Value = 0.08080657385887599
$ws.Range("T9").Value = 0.080806573858876

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("G10").Value = 1.476376
$ws.Range("H10").Value = 4.429128
$ws.Range("I10").Value = 0.1554395296814417
$ws.Range("J10").Value = 0.1554395296814417
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 5.002328333333334
$ws.Range("N10").Value = 15.006985
$ws.Range("O10").Value = 0.4685271629261051
$ws.Range("P10").Value = 0.4685271629261051
$ws.Range("Q10").Value = 7.385317495453335
$ws.Range("R10").Value = 66.46785745908001
$ws.Range("S10").Value = 0.07282764184821398
$ws.Range("T10").Value = 0.072827641848214

# Remove the now-obsolete Resolving-Mac rows (10:13 originally -> rows 11-13 removed)
$ws.Rows("11:13").Delete()
